# Remove the redundant "(vol)" suffix from organ-at-risk names in column A,
# since volume and number-of-parts are now read from the same file instead
# of from CP (the suffix was used to disambiguate the source before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 1; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value2
    if ($current -ne $null -and $current.ToString().EndsWith("(vol)")) {
        $text = $current.ToString()
        $cell.Value2 = $text.Substring(0, $text.Length - 5)
    }
}
